$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.482.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.243.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.25%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.236.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.08%  "
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.465"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.776.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.10%  "
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.236.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.469.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.729"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E33").Value = "  -3.44%  "
$ws.Range("E34").Value = "  -4.53%  "
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  -4.41%  "
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "423.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.54%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.46%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.982.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.79%  "
$ws.Range("E44").Value = "  -8.48%  "
$ws.Range("E45").Value = "  +2.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.81%  "
